# Auto-generated edit script applying literal value updates described in the diff.
# Each Leve market-data row has its H/I/J/K/L/M/N price & profit columns
# refreshed to newly observed values. Some rows gain or lose the optional
# M (LeveProfitNQ) / N (LeveProfitHQ) cell depending on whether a profit value applies.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1190.4546
$ws.Range("I18").Value = 709.5
$ws.Range("K18").Value = 709.5
$ws.Range("M18").Value = -425.5

$ws.Range("H69").Value = 3682.8572
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 3556
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 10668
$ws.Range("M69").Value = -11126
$ws.Range("N69").Value = -12416

$ws.Range("H72").Value = 3682.8572
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 3556
$ws.Range("K72").Value = 36000
$ws.Range("L72").Value = 32004
$ws.Range("M72").Value = -31632
$ws.Range("N72").Value = -40740

$ws.Range("H80").Value = 68132.53
$ws.Range("I80").Value = 1378.4
$ws.Range("J80").Value = 101509.6
$ws.Range("K80").Value = 4135.200000000001
$ws.Range("L80").Value = 304528.8
$ws.Range("M80").Value = -3137.200000000001
$ws.Range("N80").Value = -306524.8

$ws.Range("H83").Value = 68132.53
$ws.Range("I83").Value = 1378.4
$ws.Range("J83").Value = 101509.6
$ws.Range("K83").Value = 12405.6
$ws.Range("L83").Value = 913586.4
$ws.Range("M83").Value = -7413.6
$ws.Range("N83").Value = -923570.4

$ws.Range("H132").Value = 5107026
$ws.Range("I132").Value = 5957830
$ws.Range("J132").Value = 2199.4285
$ws.Range("K132").Value = 17873490
$ws.Range("L132").Value = 6598.2855
$ws.Range("M132").Value = -17870960
$ws.Range("N132").Value = -11658.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5026.7446
$ws.Range("I32").Value = 4353.43
$ws.Range("J32").Value = 19503
$ws.Range("K32").Value = 4353.43
$ws.Range("L32").Value = 19503
$ws.Range("M32").Value = -4066.43
$ws.Range("N32").Value = -20077

$ws.Range("H46").Value = 2484
$ws.Range("J46").Value = 2826
$ws.Range("L46").Value = 2826
$ws.Range("N46").Value = -3464

$ws.Range("H97").Value = 32187.938
$ws.Range("I97").Value = 44082.695
$ws.Range("J97").Value = 1790.2222
$ws.Range("K97").Value = 44082.695
$ws.Range("L97").Value = 1790.2222
$ws.Range("M97").Value = -43586.695
$ws.Range("N97").Value = -2782.2222

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 378.125
$ws.Range("I22").Value = 378.125
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 378.125
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = -205.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2116.1973
$ws.Range("I31").Value = 1402.7931
$ws.Range("J31").Value = 2608.7856
$ws.Range("K31").Value = 1402.7931
$ws.Range("L31").Value = 2608.7856
$ws.Range("M31").Value = -1107.7931
$ws.Range("N31").Value = -3198.7856

$ws.Range("H34").Value = 2116.1973
$ws.Range("I34").Value = 1402.7931
$ws.Range("J34").Value = 2608.7856
$ws.Range("K34").Value = 1402.7931
$ws.Range("L34").Value = 2608.7856
$ws.Range("M34").Value = -1200.7931
$ws.Range("N34").Value = -3012.7856

$ws.Range("H58").Value = 2245.4546
$ws.Range("I58").Value = 2275
$ws.Range("J58").Value = 2166.6667
$ws.Range("K58").Value = 2275
$ws.Range("L58").Value = 2166.6667
$ws.Range("M58").Value = -2072
$ws.Range("N58").Value = -2572.6667

$ws.Range("H86").Value = 3833.75
$ws.Range("I86").Value = 3500
$ws.Range("J86").Value = 4072.1428
$ws.Range("K86").Value = 3500
$ws.Range("L86").Value = 4072.1428
$ws.Range("M86").Value = -2377
$ws.Range("N86").Value = -6318.1428

$ws.Range("H89").Value = 3833.75
$ws.Range("I89").Value = 3500
$ws.Range("J89").Value = 4072.1428
$ws.Range("K89").Value = 17500
$ws.Range("L89").Value = 20360.714
$ws.Range("M89").Value = -11884
$ws.Range("N89").Value = -31592.714

$ws.Range("H136").Value = 2245.4546
$ws.Range("I136").Value = 2275
$ws.Range("J136").Value = 2166.6667
$ws.Range("K136").Value = 6825
$ws.Range("L136").Value = 6500.000100000001
$ws.Range("M136").Value = -4275
$ws.Range("N136").Value = -11600.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 641
$ws.Range("J61").Value = 641
$ws.Range("L61").Value = 1923
$ws.Range("N61").Value = -2353

$ws.Range("H97").Value = 1900.6666
$ws.Range("I97").Value = 1625
$ws.Range("J97").Value = 2452
$ws.Range("K97").Value = 4875
$ws.Range("L97").Value = 7356
$ws.Range("M97").Value = -4379
$ws.Range("N97").Value = -8348

$ws.Range("H137").Value = 39113.445
$ws.Range("I137").Value = 60298.234
$ws.Range("J137").Value = 3099.3
$ws.Range("K137").Value = 180894.702
$ws.Range("L137").Value = 9297.900000000001
$ws.Range("M137").Value = -175794.702
$ws.Range("N137").Value = -19497.9

$ws.Range("H140").Value = 4709.129
$ws.Range("I140").Value = 5706.227
$ws.Range("J140").Value = 2271.7778
$ws.Range("K140").Value = 17118.681
$ws.Range("L140").Value = 6815.3334
$ws.Range("M140").Value = -11938.681
$ws.Range("N140").Value = -17175.3334

$ws.Range("H141").Value = 18837.334
$ws.Range("I141").Value = 21404.8
$ws.Range("K141").Value = 64214.39999999999
$ws.Range("M141").Value = -59034.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1555.5385
$ws.Range("I113").Value = 1162.5
$ws.Range("J113").Value = 1730.2222
$ws.Range("K113").Value = 1162.5
$ws.Range("L113").Value = 1730.2222
$ws.Range("M113").Value = 1007.5
$ws.Range("N113").Value = -6070.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4039.4
$ws.Range("I7").Value = 2734.6667
$ws.Range("J7").Value = 4598.5713
$ws.Range("K7").Value = 2734.6667
$ws.Range("L7").Value = 4598.5713
$ws.Range("M7").Value = -2622.6667
$ws.Range("N7").Value = -4822.5713

$ws.Range("H46").Value = 327401.28
$ws.Range("I46").Value = 1193.3334
$ws.Range("J46").Value = 779073.9
$ws.Range("K46").Value = 1193.3334
$ws.Range("L46").Value = 779073.9
$ws.Range("M46").Value = -1005.3334
$ws.Range("N46").Value = -779449.9

$ws.Range("H92").Value = 29500
$ws.Range("J92").Value = 29500
$ws.Range("L92").Value = 29500
$ws.Range("N92").Value = -34492

$ws.Range("H93").Value = 1975.8948
$ws.Range("I93").Value = 1881.2
$ws.Range("K93").Value = 1881.2
$ws.Range("M93").Value = -633.2

$ws.Range("H98").Value = 23333.334
$ws.Range("J98").Value = 23333.334
$ws.Range("L98").Value = 23333.334
$ws.Range("N98").Value = -29323.334

$ws.Range("H99").Value = 28600
$ws.Range("I99").Value = 25900
$ws.Range("K99").Value = 25900
$ws.Range("M99").Value = -22905

$ws.Range("H100").Value = 2414.4443
$ws.Range("I100").Value = 2183.3333
$ws.Range("J100").Value = 2876.6667
$ws.Range("K100").Value = 2183.3333
$ws.Range("L100").Value = 2876.6667
$ws.Range("M100").Value = -1642.3333
$ws.Range("N100").Value = -3958.6667

$ws.Range("H101").Value = 15024.889
$ws.Range("J101").Value = 15024.889
$ws.Range("L101").Value = 15024.889
$ws.Range("N101").Value = -21514.889

$ws.Range("H102").Value = 43779
$ws.Range("J102").Value = 43779
$ws.Range("L102").Value = 43779
$ws.Range("N102").Value = -50269

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = $null
$ws.Range("N103").Value = 0

$ws.Range("H105").Value = 45000
$ws.Range("J105").Value = 45000
$ws.Range("L105").Value = 45000
$ws.Range("N105").Value = -51988

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = $null
$ws.Range("N106").Value = 0

$ws.Range("H122").Value = 2672.1333
$ws.Range("I122").Value = 2621.6924
$ws.Range("K122").Value = 7865.0772
$ws.Range("M122").Value = -5415.0772

$ws.Range("H126").Value = 4039.4
$ws.Range("I126").Value = 2734.6667
$ws.Range("J126").Value = 4598.5713
$ws.Range("K126").Value = 8204.000100000001
$ws.Range("L126").Value = 13795.7139
$ws.Range("M126").Value = -5734.000100000001
$ws.Range("N126").Value = -18735.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 38463036
$ws.Range("I62").Value = 76923070
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 76923070
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -76922446
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 38463036
$ws.Range("I65").Value = 76923070
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 384615350
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -384612230
$ws.Range("N65").Value = -21240

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = $null
$ws.Range("N92").Value = 0

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = $null
$ws.Range("N93").Value = 0

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = $null
$ws.Range("N95").Value = 0

$ws.Range("H96").Value = 125001590
$ws.Range("I96").Value = 200001520
$ws.Range("J96").Value = 1726.6666
$ws.Range("K96").Value = 200001520
$ws.Range("L96").Value = 1726.6666
$ws.Range("M96").Value = -200000147
$ws.Range("N96").Value = -4472.6666

$ws.Range("H99").Value = 34900
$ws.Range("J99").Value = 34900
$ws.Range("L99").Value = 34900
$ws.Range("N99").Value = -40890

$ws.Range("H100").Value = 250960
$ws.Range("I100").Value = 333446.66
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 666893.3199999999
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = -666352.3199999999
$ws.Range("N100").Value = -8082

$ws.Range("H101").Value = 16500
$ws.Range("J101").Value = 16500
$ws.Range("L101").Value = 16500
$ws.Range("N101").Value = -22990

$ws.Range("H102").Value = 39990
$ws.Range("J102").Value = 39990
$ws.Range("L102").Value = 39990
$ws.Range("N102").Value = -46480

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = $null
$ws.Range("N103").Value = 0

$ws.Range("H104").Value = 16500
$ws.Range("J104").Value = 16500
$ws.Range("L104").Value = 16500
$ws.Range("N104").Value = -23488

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = $null
$ws.Range("N105").Value = 0
